# Update "Förändrad" (Changed) date column (C) for rows 2-8 on the active sheet.
# Old serial date value 45212 (2023-10-13) -> new serial date value 45221 (2023-10-22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
